# Saldo.xlsx update: remove several low-priority account rows, correct
# ANILSON's balance, and re-insert the CINTIA account lower in the sheet
# with an updated balance.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# --- Correct ANILSON's balance (account 004385806), while rows are still
#     at their original (pre-delete) positions: row 10 ---
$ws.Range("C10").Value = 6670.46

# --- Delete rows (original row numbers, highest first so the row number
#     of any row not yet processed is never shifted) ---
$rowsToDelete = @(15, 14, 13, 12, 11, 9, 8, 7, 4)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

# --- Re-insert account 004927044 / CINTIA right after the DIEGO /
#     005654767 row (originally row 40, now row 40-9=31 after the 9
#     deletions above), using its new, much lower balance ---
$ws.Rows(32).Insert()
$ws.Cells.Item(32, 1).NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "004927044"
$ws.Cells.Item(32, 2).Value = "CINTIA"
$ws.Cells.Item(32, 3).Value = 316.61
